$d = $word.ActiveDocument

# Locate the three reference paragraphs in the "peer-reviewed publicaties" list
# (searched by stable text fragments rather than hard-coded indices, so the
# script is resilient to minor renumbering).
$oldArentsenPara = $null
$duplicateArentsenPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Contains("Energy policy and nuclear power - 20 years after the Chernobyl disaster")) {
        $oldArentsenPara = $p
    }
    elseif ($t.Contains("095830506778119407")) {
        $duplicateArentsenPara = $p
    }
}

$fullCitation = "Arentsen, M. J. (2006). CONTESTED TECHNOLOGY: Nuclear Power in the Netherlands. Energy & Environment (Essex, England), 17(3), 373" + [char]0x2013 + "382. https://doi.org/10.1260/095830506778119407"

# Remove the now-redundant trailing duplicate paragraph along with the blank
# paragraph that follows it (so only a single blank line remains between
# "van Leeuwen..." and "Zijlstra..."). Deleting forward into the following
# blank paragraph's own mark (rather than swallowing the preceding blank's
# mark) is what actually merges the paragraphs in this runtime.
$followingBlankIndex = $duplicateArentsenPara.Index + 1
$followingBlank = $d.Paragraphs($followingBlankIndex)
$deleteRange = $d.Range($duplicateArentsenPara.Range.Start, $followingBlank.Range.End)
$deleteRange.Delete()

# Replace the messy, multi-run "ARENTSEN..." paragraph's content with the
# clean, complete citation (collapses to a single run). Done after the
# deletion above so this still-valid Range (it sits before the deleted
# text) is rewritten into the single final run.
$replaceRange = $d.Range($oldArentsenPara.Range.Start, $oldArentsenPara.Range.End - 1)
$replaceRange.Text = $fullCitation

Write-Output "done"
